$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-5 need to be cyclically shifted down by one, with row 5's
# data wrapping around to become the new row 2 (row2->row3, row3->row4,
# row4->row5, row5->row2). Use row 6 as scratch space, clearing destination
# ranges before each copy so that blank/absent source cells properly clear
# any stale values left behind in the destination.

$ws.Range("A5:AY5").Copy($ws.Range("A6:AY6"))

$ws.Range("A5:AY5").Clear()
$ws.Range("A4:AY4").Copy($ws.Range("A5:AY5"))

$ws.Range("A4:AY4").Clear()
$ws.Range("A3:AY3").Copy($ws.Range("A4:AY4"))

$ws.Range("A3:AY3").Clear()
$ws.Range("A2:AY2").Copy($ws.Range("A3:AY3"))

$ws.Range("A2:AY2").Clear()
$ws.Range("A6:AY6").Copy($ws.Range("A2:AY2"))

$ws.Range("A6:AY6").Clear()
